$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Hello again! Ukuba umzali kungaba yingcindezi, futhi iskhathi sokuthi uzinakekele nawe!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Sawubona futhi! Ukuba umzali kungaba yingcindezi, futhi iskhathi sokuthi uzinakekele nawe!",
    2) | Out-Null

$d.Content.Find.Execute(
    "Here is a simple stretching and movement activity that may help you with stress.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Nansi ndlela yokuzilula nokunyakaza engakusiza ngengcindezi onayo.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Stand up and stretch your arms up to the sky.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Sukuma ulule zandla zakho ziyephezulu.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Stretch to both sides.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Zilule nhlangothi zombili.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Stretch to the front and to the back.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Zilulele ngaphimbili nange muva.",
    2) | Out-Null
